$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.815.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "'1.889.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'239.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.4763"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.75%  "

$ws.Range("D8").Value = "'0.2881"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.03%  "

$ws.Range("D9").Value = "'0.06602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.58%  "

$ws.Range("D10").Value = "'18.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.66%  "

$ws.Range("D11").Value = "'99.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.21%  "

$ws.Range("D12").Value = "'1.880.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").Value = "'0.07607"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "'5.143"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.00%  "

$ws.Range("D15").Value = "'0.6627"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.33%  "

$ws.Range("D16").Value = "'309.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +35.47%  "

$ws.Range("D17").Value = "'30.816.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").Value = "'13.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.63%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").Value = "'0.000007588"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.89%  "

$ws.Range("D21").Value = "'2.123.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'5.130"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").Value = "'6.195"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.81%  "

$ws.Range("D25").Value = "'9.318"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("D26").Value = "'167.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").Value = "'20.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.99%  "

$ws.Range("D28").Value = "'1.950"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.72%  "

$ws.Range("D29").Value = "'0.1077"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.69%  "

$ws.Range("D30").Value = "'1.359"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "'4.187"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.44%  "

$ws.Range("D32").Value = "'3.994"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.56%  "

$ws.Range("D33").Value = "'0.05061"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.54%  "

$ws.Range("D34").Value = "'1.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.24%  "

$ws.Range("D35").Value = "'0.7316"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.50%  "

$ws.Range("D36").Value = "'2.716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").Value = "'0.01956"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.48%  "

$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("D39").Value = "'2.081"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.00%  "

$ws.Range("D40").Value = "'0.9047"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("D41").Value = "'107.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").Value = "'0.9998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "'0.4220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.22%  "

$ws.Range("D44").Value = "'5.657"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "

$ws.Range("D45").Value = "'7.415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.46%  "

$ws.Range("D46").Value = "'65.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.53%  "

$ws.Range("D47").Value = "'9.071"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.86%  "

$ws.Range("D48").Value = "'0.1228"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("D49").Value = "'34.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.05%  "

$ws.Range("D50").Value = "'0.05631"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").Value = "'1.393"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
